$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Outcomes")

# Insert a new row before the existing row 67 ("Social worker ethnicity"),
# shifting it (and nothing else, since it was the last data row) down to row 68.
$ws.Rows.Item(67).Insert()

# Populate the newly inserted row 67 with the "Social worker stability" metric
# used by the LA chart / Stats neighbours chart module.
$ws.Range("A67").Value = 66
$ws.Range("B67").Value = "Enablers"
$ws.Range("C67").Value = "Enabler: The workforce is equipped and effective"
$ws.Range("D67").Value = "Quality of support for children and families"
$ws.Range("E67").Value = "Social worker stability"
$ws.Range("F67").Value = "sw_stability"
$ws.Range("G67").Value = "percent"
$ws.Range("H67").Value = "Percent"
$ws.Range("I67").Value = "list('cla_group'='CLA on 31 March','sw_stability'='3 or more social workers during the year')"
$ws.Range("J67").Value = "list()"

# The row that used to be #67 (sort_order 66) is now row 68; bump its
# sort_order value to 67 to keep the sequential numbering intact.
$ws.Range("A68").Value = 67

# Update the view state of the Outcomes sheet to match where the author
# left off after adding the new row.
$ws.Activate()
$win = $excel.ActiveWindow
$win.Zoom = 130
$win.ScrollRow = 37
$win.ScrollColumn = 4
$ws.Range("E67").Select()
